$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.485.83'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.92%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.918.11'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.50%  '

$ws.Range("E4").Value = '  -0.14%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '326.32'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.77%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.18%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4752'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.69%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4095'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.53%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '47.94'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.75%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08056'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.84%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.012'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.42%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.51'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.32%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.869.82'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.06%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.935'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.28%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.160'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.33%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '89.61'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.50%  '

$ws.Range("E17").Value = '  -0.17%  '

$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001033'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.29%  '

$ws.Range("B19").Value = 'TRON'
$ws.Range("C19").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06597'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.30%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.76'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.38%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.000'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.18%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '29.497.64'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.81%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.541'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.54%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.50'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.28%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.206'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.54%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.133.37'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.87%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '154.68'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.78%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.83'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.64%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.951'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +9.61%  '

$ws.Range("E30").Value = '  +0.35%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '117.63'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.23%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.059'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +7.89%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09552'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.82%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.433'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.93%  '

$ws.Range("E35").Value = '  -1.00%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.409'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.37%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06121'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.82%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02260'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.99%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.322'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.27%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.171'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.02%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5893'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.92%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.557'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +10.83%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1846'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.34%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.16'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.05%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.08046'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +14.69%  '

$ws.Range("E46").Value = '  +2.78%  '

$ws.Range("B47").Value = 'Decentraland'
$ws.Range("C47").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5556'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.13%  '

$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '12.13'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.96%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.937'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.63%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '113.20'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.84%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '44.77'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.29%  '
